$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# GPLIM-2588 Fix spreadsheet headers.
$ws.Cells.Item(1, 1).Value2 = "Specimen_Number"
$ws.Cells.Item(1, 6).Value2 = "SAMPLE_TYPE"

# Re-apply A1's original formatting (quote-prefixed header style), since
# writing a new value to the cell would otherwise re-derive its style.
$ws.Cells.Item(1, 2).Copy() | Out-Null
$ws.Cells.Item(1, 1).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A2").Select()
